# Update "想去人数" (interested-attendee count) values in column F
# across all four sheets, reflecting newer scrape data.

$wb = $excel.ActiveWorkbook

# --- Sheet 1: 展览 ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value  = 135
$ws1.Range("F5").Value  = 1954
$ws1.Range("F7").Value  = 4054
$ws1.Range("F13").Value = 97
$ws1.Range("F14").Value = 2175
$ws1.Range("F15").Value = 391
$ws1.Range("F16").Value = 650315
$ws1.Range("F18").Value = 484
$ws1.Range("F21").Value = 539
$ws1.Range("F23").Value = 2172
$ws1.Range("F25").Value = 2670
$ws1.Range("F26").Value = 1537
$ws1.Range("F27").Value = 774
$ws1.Range("F28").Value = 1515
$ws1.Range("F30").Value = 1076
$ws1.Range("F34").Value = 1338
$ws1.Range("F36").Value = 2123
$ws1.Range("F41").Value = 2551

# --- Sheet 2: 演出 ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F11").Value = 144407
$ws2.Range("F12").Value = 144407
$ws2.Range("F16").Value = 23
$ws2.Range("F27").Value = 528
$ws2.Range("F32").Value = 323
$ws2.Range("F43").Value = 7

# --- Sheet 3: 本地生活 ---
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F5").Value = 237
$ws3.Range("F8").Value = 1153

# --- Sheet 4: 全部类型 ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F7").Value  = 135
$ws4.Range("F9").Value  = 4054
$ws4.Range("F15").Value = 2175
$ws4.Range("F17").Value = 391
$ws4.Range("F18").Value = 650316
$ws4.Range("F22").Value = 144407
$ws4.Range("F23").Value = 484
$ws4.Range("F26").Value = 539
$ws4.Range("F28").Value = 2172
$ws4.Range("F30").Value = 2670
$ws4.Range("F31").Value = 1537
$ws4.Range("F33").Value = 1515
$ws4.Range("F37").Value = 1076
$ws4.Range("F42").Value = 1338
$ws4.Range("F44").Value = 2123
$ws4.Range("F46").Value = 323
$ws4.Range("F47").Value = 323
$ws4.Range("F48").Value = 2551
